# add_save_column.ps1
# Adds a new "Save" column (H) to the s_vals sheet, mirroring the
# header formatting already used by the other header cells (B1:G1)
# and writing a numeric value for the single data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1, which carries
# the bold/border/center style used by every header) onto the new H1
# header cell, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cell for row 2 — numeric flag column.
$ws.Range("H2").Value = 1
